$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update intro text date (04 April 2025 -> 11 April 2025)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 11 April 2025"

# 2. Fix statuses for two existing rows (provisional -> confirmed)
$ws.Range("D8").Value = "confirmed"
$ws.Range("D10").Value = "confirmed"

# 3. Insert a new row at row 12 for the new publication,
#    pushing "Coroners statistics 2024" (and everything below) down by one row
$ws.Rows("12:12").Insert()
$ws.Range("A12").Value = "05 May 2025"
$ws.Range("B12").Value = "Unpaid Work Management Information, update to December 2024"
$ws.Range("C12").Value = "8 May 2025"
$ws.Range("D12").Value = "provisional"
$ws.Range("E12").Value = 19
$ws.Range("F12").Value = "standard"

# 4. Extend conditional formatting ranges to cover the new last row (64)
$fcs = $ws.Range("A5").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A5:F64"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("A5:F64"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("A5:F64"))
$fcs.Item(4).ModifyAppliesToRange($ws.Range("A5:A64"))
$fcs.Item(5).ModifyAppliesToRange($ws.Range("A5:A64"))
